$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "282.51"
$ws.Range("D2").Style = "Normal"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "10"
$ws.Range("G2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "20.88"
$ws.Range("D3").Style = "Normal"

$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "10"
$ws.Range("G3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.246"
$ws.Range("D4").Style = "Normal"

$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "10"
$ws.Range("G4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06180"
$ws.Range("D5").Style = "Normal"

$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "10"
$ws.Range("G5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.586"
$ws.Range("D6").Style = "Normal"

$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "10"
$ws.Range("G6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.562"
$ws.Range("D7").Style = "Normal"

$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "10"
$ws.Range("G7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.492"
$ws.Range("D8").Style = "Normal"

$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "10"
$ws.Range("G8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8158"
$ws.Range("D9").Style = "Normal"

$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "10"
$ws.Range("G9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01390"
$ws.Range("D10").Style = "Normal"

$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "10"
$ws.Range("G10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1638"
$ws.Range("D11").Style = "Normal"

$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "10"
$ws.Range("G11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08320"
$ws.Range("D12").Style = "Normal"

$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "10"
$ws.Range("G12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03556"
$ws.Range("D13").Style = "Normal"

$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "10"
$ws.Range("G13").Style = "Normal"

$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "10"
$ws.Range("G14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09137"
$ws.Range("D15").Style = "Normal"

$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "10"
$ws.Range("G15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.724"
$ws.Range("D16").Style = "Normal"

$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "10"
$ws.Range("G16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001640"
$ws.Range("D17").Style = "Normal"

$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "10"
$ws.Range("G17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04689"
$ws.Range("D18").Style = "Normal"

$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "10"
$ws.Range("G18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006420"
$ws.Range("D19").Style = "Normal"

$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "10"
$ws.Range("G19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.006189"
$ws.Range("D20").Style = "Normal"

$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "10"
$ws.Range("G20").Style = "Normal"

$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "10"
$ws.Range("G21").Style = "Normal"

$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "10"
$ws.Range("G22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.821"
$ws.Range("D23").Style = "Normal"

$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "10"
$ws.Range("G23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.338"
$ws.Range("D24").Style = "Normal"

$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "10"
$ws.Range("G24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3377"
$ws.Range("D25").Style = "Normal"

$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "10"
$ws.Range("G25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1249"
$ws.Range("D26").Style = "Normal"

$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "10"
$ws.Range("G26").Style = "Normal"

$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "10"
$ws.Range("G27").Style = "Normal"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "10"
$ws.Range("G28").Style = "Normal"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "10"
$ws.Range("G29").Style = "Normal"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "10"
$ws.Range("G30").Style = "Normal"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "10"
$ws.Range("G31").Style = "Normal"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "10"
$ws.Range("G32").Style = "Normal"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "10"
$ws.Range("G33").Style = "Normal"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "10"
$ws.Range("G34").Style = "Normal"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "10"
$ws.Range("G35").Style = "Normal"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "10"
$ws.Range("G36").Style = "Normal"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "10"
$ws.Range("G37").Style = "Normal"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "10"
$ws.Range("G38").Style = "Normal"

$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "10"
$ws.Range("G39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04684"
$ws.Range("D40").Style = "Normal"

$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "10"
$ws.Range("G40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007127"
$ws.Range("D41").Style = "Normal"

$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "10"
$ws.Range("G41").Style = "Normal"

$ws.Range("B42").Value = "CEJI"

$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004398"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "10"
$ws.Range("G42").Style = "Normal"

$ws.Range("B43").Value = "BKEXToken"

$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1101"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "42BKEXTokenBKK"

$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "10"
$ws.Range("G43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01111"
$ws.Range("D44").Style = "Normal"

$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "10"
$ws.Range("G44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006323"
$ws.Range("D45").Style = "Normal"

$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "10"
$ws.Range("G45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D46").Style = "Normal"

$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "10"
$ws.Range("G46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9995"
$ws.Range("D47").Style = "Normal"

$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "10"
$ws.Range("G47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002924"
$ws.Range("D48").Style = "Normal"

$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "10"
$ws.Range("G48").Style = "Normal"

$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "10"
$ws.Range("G49").Style = "Normal"

$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "10"
$ws.Range("G50").Style = "Normal"

$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "10"
$ws.Range("G51").Style = "Normal"
